# Apply updated odds values to row 2 of the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.45
$ws.Range("I2").Value = 3.1
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 3.6
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 23
$ws.Range("AM2").Value = 34
$ws.Range("AO2").Value = 15
$ws.Range("AQ2").Value = 51
$ws.Range("AW2").Value = 4.75
